$wb = $excel.ActiveWorkbook

# --- Sheet 1: ClientInfo ---
$ws1 = $wb.Worksheets.Item("ClientInfo")

# Row 2
$ws1.Range("A2").Value = 750
$ws1.Range("B2").Value = "MIKHAIL"
$ws1.Range("C2").Value = "NOVIKAU"
$ws1.Range("D2").Value = 28
$ws1.Range("E2").Value = 30001

# Row 3
$ws1.Range("A3").Value = 1449
$ws1.Range("B3").Value = "TEST"
$ws1.Range("C3").Value = "TEST"
$ws1.Range("D3").Value = 25
$ws1.Range("E3").Value = 1000

# Row 4 no longer present - delete it entirely
$ws1.Rows.Item(4).Delete()

# --- Sheet 2: ClientAccountInfo ---
$ws2 = $wb.Worksheets.Item("ClientAccountInfo")

# Row 2
$ws2.Range("B2").Value = "MELMENX"
$ws2.Range("C2").Value = "2127979Z"

# Row 3
$ws2.Range("B3").Value = "TEST"
$ws2.Range("C3").Value = "TEST"

# Row 4 no longer present - delete it entirely
$ws2.Rows.Item(4).Delete()
